$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Write the new cell values in the exact order needed so that the shared
#    string table is built up in the same sequence as the target workbook.
# ---------------------------------------------------------------------------

# Row 2 additions (D2, E2)
$ws.Range("D2").Value = "CalvinDerekAG"
$ws.Range("E2").Value = "Calvin Derek"

# Row 1 additions (D1, E1) - header row
$ws.Range("D1").Value = "ApprovalGroupName"
$ws.Range("E1").Value = "ApprovalManager"

# Row 3 (new data row)
$ws.Range("A3").Value = "AddApprovalGroupPO_GuideWellGroupInc"
$ws.Range("B3").Value = 3670
$ws.Range("C3").Value = "GuideWell Group Inc."
$ws.Range("D3").Value = "Primary Owner"
$ws.Range("E3").Value = "Dunston"

# Row 4 (new data row)
$ws.Range("A4").Value = "AddApprovalGroupPO_SEIT0019"
$ws.Range("B4").Value = 3858
$ws.Range("C4").Value = "Burns & McDonnell"
$ws.Range("D4").Value = "Primary Owner"
$ws.Range("E4").Value = "Acton Sara"

# Row 5 (new data row)
$ws.Range("A5").Value = "AddAutomationGroup_WyndhamWorldwide"
$ws.Range("B5").Value = 392
$ws.Range("C5").Value = "Wyndham Destinations"
$ws.Range("D5").Value = "Automation Group"
$ws.Range("E5").Value = "Perez Madalene"

# ---------------------------------------------------------------------------
# 2. Apply the existing bordered style (same as columns A-C) to the new
#    D/E columns. We copy the format from a cell that already carries that
#    style so the engine reuses the same cellXfs/border entry instead of
#    creating a new one.
# ---------------------------------------------------------------------------

# Body rows 2-14 (D:E) -> reuse the plain bordered style already used by A:C
$ws.Range("A2").Copy()
$ws.Range("D2:E14").PasteSpecial(-4122)

# Re-apply the values (PasteSpecial(formats) should not disturb them, but do
# it defensively in case paste touched D2/E2 which already had values).
$ws.Range("D2").Value = "CalvinDerekAG"
$ws.Range("E2").Value = "Calvin Derek"
$ws.Range("D3").Value = "Primary Owner"
$ws.Range("E3").Value = "Dunston"
$ws.Range("D4").Value = "Primary Owner"
$ws.Range("E4").Value = "Acton Sara"
$ws.Range("D5").Value = "Automation Group"
$ws.Range("E5").Value = "Perez Madalene"

# Header cells D1:E1 -> bordered style plus an explicit (no-op) fill flag,
# mirroring the workbook's third cellXfs entry.
$ws.Range("A2").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = "ApprovalGroupName"
$ws.Range("E1").Value = "ApprovalManager"
$ws.Range("D1:E1").Interior.ColorIndex = -4142

# ---------------------------------------------------------------------------
# 3. Resize columns A, C, D, E to fit their (now longer) contents.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()

# ---------------------------------------------------------------------------
# 4. Update the active selection to match the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("A5").Select()
